{"js": "// \"Add intro text and summary chart\" \u2014 template style updates:\n//   1. Heading1 paragraph style: force a page break before each Heading 1\n//      (so the new intro/summary sections start on their own page).\n//   2. Heading1Char (the linked run style for Heading 1): switch the\n//      heading font from \"Source Sans Pro\" to \"Arial\".\n//   3. TableCharttitle (used above the new summary chart/table): center\n//      the chart title text.\n\nconst styles = context.document.getStyles();\n\n// 1. Heading1 \u2014 add pageBreakBefore to the paragraph properties.\nconst heading1 = styles.getByName(\"Heading1\");\nheading1.paragraphFormat.pageBreakBefore = true;\n\n// 2. Heading1Char \u2014 swap the heading font to Arial (ascii + hAnsi).\nconst heading1Char = styles.getByName(\"Heading1Char\");\nheading1Char.font.name = \"Arial\";\n\n// 3. TableChart title \u2014 center-align the chart title paragraph.\nconst tableChartTitle = styles.getByName(\"TableCharttitle\");\ntableChartTitle.paragraphFormat.alignment = Word.Alignment.centered;\n\nawait context.sync();\n", "ps1": "# \"Add intro text and summary chart\" \u2014 template style updates:\n#   1. Heading1 paragraph style: force a page break before each Heading 1\n#      (so the new intro/summary sections start on their own page).\n#   2. Heading1Char (the linked run style for Heading 1): switch the\n#      heading font from \"Source Sans Pro\" to \"Arial\".\n#   3. TableCharttitle (used above the new summary chart/table): center\n#      the chart title text.\n\n$d = $word.ActiveDocument\n\n# 1. Heading1 - add pageBreakBefore to the paragraph properties.\n$heading1 = $d.Styles.Item(\"Heading1\")\n$heading1.ParagraphFormat.PageBreakBefore = $true\n\n# 2. Heading1Char - swap the heading font to Arial (ascii + hAnsi).\n$heading1Char = $d.Styles.Item(\"Heading1Char\")\n$heading1Char.Font.Name = \"Arial\"\n\n# 3. TableChart title - center-align the chart title paragraph.\n$tableChartTitle = $d.Styles.Item(\"TableCharttitle\")\n$tableChartTitle.ParagraphFormat.Alignment = \"wdAlignParagraphCenter\"\n"}
